$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "DESeq2/LimmaVoom/EdgeR"
$ws.Range("B1").Value = "DESeq2/EdgeR"
$ws.Range("C1").Value = "DESeq2/LimmaVoom"
$ws.Range("D1").Value = "EdgeR/LimmaVoom"
$ws.Range("E1").Value = "DESeq2"
$ws.Range("F1").Value = "LimmaVoom"
$ws.Range("G1").Value = "EdgeR"
